$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30-47 down to 31-48.
# Excel copies the formatting of the row above on insert, matching the
# original file's per-row style (date column D uses style index 2).
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with the new weekly price record.
$ws.Range("A30").Value = 4
$ws.Range("B30").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C30").Value = "Los Lagos"
$ws.Range("D30").Value = 44827
$ws.Range("E30").Value = 10
$ws.Range("F30").Value = 100112012
$ws.Range("G30").Value = "Espinaca"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 35
$ws.Range("K30").Value = 13000
$ws.Range("L30").Value = 13000
$ws.Range("M30").Value = 13000
$ws.Range("N30").Value = "$/cuna 10 kilos"
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 1300
$ws.Range("Q30").Value = 10
$ws.Range("R30").Value = "Hortaliza"
